# Correcting problems with name changes:
# Renumber column C for rows 28-53 from 27..52 down to 1..26,
# and update the active selection to reflect the new working range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 28; $row -le 53; $row++) {
    $ws.Cells.Item($row, 3).Value = $row - 27
}

$ws.Range("C28:C53").Select()
